$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data generator folded in one more source JSON for the
# time-bucket analysis. As a result, the two existing articles that are
# closest to the event swap position in the table: the "day_0" article
# ("Noblesville school shooting: New details about the day of the shooting")
# now sits in row 2 and the "day_2_to_30" article ("Indiana law to be
# reviewed after Noblesville shooting suspect tried as a child") now sits in
# row 3. Row 4 ("Indiana middle schooler...") is untouched.

# Capture the current (pre-edit) values of row 2 and row 3, columns A-D.
$a2 = $ws.Range("A2").Value()
$b2 = $ws.Range("B2").Value()
$c2 = $ws.Range("C2").Value()
$d2 = $ws.Range("D2").Value()
$e2 = $ws.Range("E2").Value()

$a3 = $ws.Range("A3").Value()
$b3 = $ws.Range("B3").Value()
$c3 = $ws.Range("C3").Value()
$d3 = $ws.Range("D3").Value()
$e3 = $ws.Range("E3").Value()

$e4 = $ws.Range("E4").Value()

# Write the swapped A-D values back (row 2 gets row 3's data and vice versa).
$ws.Range("A2").Value = $a3
$ws.Range("B2").Value = $b3
$ws.Range("C2").Value = $c3
$ws.Range("D2").Value = $d3

$ws.Range("A3").Value = $a2
$ws.Range("B3").Value = $b2
$ws.Range("C3").Value = $c2
$ws.Range("D3").Value = $d2

# Rebuild the hyperlinks for column E so each row's link follows its text
# (E2 now points to the indystar.com article, E3 to the theindychannel.com
# article; E4/abcnews.go.com is re-added unchanged).
$hls = $ws.Hyperlinks
$hls.Delete()
$hls.Add($ws.Range("E2"), $e3) | Out-Null
$hls.Add($ws.Range("E3"), $e2) | Out-Null
$hls.Add($ws.Range("E4"), $e4) | Out-Null

# The Add() call above only rewrites the hyperlink relationship/target; the
# cell's own displayed text needs to be set explicitly so it matches.
$ws.Range("E2").Value = $e3
$ws.Range("E3").Value = $e2
$ws.Range("E4").Value = $e4

# Restore the built-in Hyperlink cell style that was lost when the
# hyperlinks collection was rebuilt.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"

Write-Output "Row 2 and row 3 swapped successfully"
